{"js": "// Update the \"Data de envio\" timestamp and reset a handful of briefing\n// answer fields back to \"N\u00e3o informado\" (exactly as described by the diff).\n//\n// The document stores each answer as: <label>\\v<value>  (the \\v == 0x0B\n// is the <w:br/> that separates the bold label run from the plain-text\n// value run). We locate the old literal value text anywhere in the body\n// and replace it in place with the new value, which preserves every\n// other run/paragraph property untouched.\n\nconst body = context.document.body;\n\n// 1) Timestamp line - unique text in the document.\nconst tsHits = body.search(\"Data de envio: 23/06/2025, 19:37:59\", { matchCase: true });\ntsHits.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < tsHits.items.length; i++) {\n  tsHits.items[i].insertText(\"Data de envio: 23/06/2025, 19:38:05\", \"Replace\");\n}\nawait context.sync();\n\n// 2) Answer fields reset to \"N\u00e3o informado\". Some old values repeat across\n// more than one field (\" Teste\" / \" Testando \") - search+replace all of\n// them is correct in every case since they all become the same new value.\nconst valueReplacements = [\n  [\" Teste\", \" N\u00e3o informado\"],\n  [\" Ser\u00e1 que funciona?\", \" N\u00e3o informado\"],\n  [\" Testando \", \" N\u00e3o informado\"],\n  [\" N\u00e3o possui logo\", \" N\u00e3o informado\"],\n  [\" Hahshss\", \" N\u00e3o informado\"],\n  [\" Meme\", \" N\u00e3o informado\"],\n  [\" Sim\", \" N\u00e3o informado\"],\n  [\" Melhorar reputa\u00e7\u00e3o online\", \" N\u00e3o informado\"],\n  [\" Aaaa\", \" N\u00e3o informado\"],\n];\n\nfor (const [oldText, newText] of valueReplacements) {\n  const hits = body.search(oldText, { matchCase: true });\n  hits.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < hits.items.length; i++) {\n    hits.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the \"Data de envio\" timestamp and reset a handful of briefing\n# answer fields back to \"N\u00e3o informado\" (exactly as described by the diff).\n#\n# Each answer is stored as two runs inside one paragraph:\n#   <b>Label:</b><br/> <value>\n# so a literal Find/Replace on the old value text (scoped to the whole\n# document body) updates exactly the right run without touching the bold\n# label run or any paragraph/run formatting.\n\n$d = $word.ActiveDocument\n\nfunction Replace-AllText($findText, $replaceText) {\n  $rng = $d.Content\n  $rng.Find.Execute(\n    [ref]$findText,    # FindText\n    [ref]$true,        # MatchCase\n    [ref]$false,       # MatchWholeWord\n    [ref]$false,       # MatchWildcards\n    [ref]$false,       # MatchSoundsLike\n    [ref]$false,       # MatchAllWordForms\n    [ref]$true,        # Forward\n    [ref]1,            # Wrap (wdFindContinue)\n    [ref]$false,       # Format\n    [ref]$replaceText, # ReplaceWith\n    [ref]2             # Replace (wdReplaceAll)\n  ) | Out-Null\n}\n\n# 1) Timestamp line (unique text in the document).\nReplace-AllText \"Data de envio: 23/06/2025, 19:37:59\" \"Data de envio: 23/06/2025, 19:38:05\"\n\n# 2) Answer fields reset to \"N\u00e3o informado\". A couple of old values repeat\n# across more than one field (\" Teste\" / \" Testando \") - replacing every\n# occurrence is correct since they all become the same new value.\nReplace-AllText \" Teste\" \" N\u00e3o informado\"\nReplace-AllText \" Ser\u00e1 que funciona?\" \" N\u00e3o informado\"\nReplace-AllText \" Testando \" \" N\u00e3o informado\"\nReplace-AllText \" N\u00e3o possui logo\" \" N\u00e3o informado\"\nReplace-AllText \" Hahshss\" \" N\u00e3o informado\"\nReplace-AllText \" Meme\" \" N\u00e3o informado\"\nReplace-AllText \" Sim\" \" N\u00e3o informado\"\nReplace-AllText \" Melhorar reputa\u00e7\u00e3o online\" \" N\u00e3o informado\"\nReplace-AllText \" Aaaa\" \" N\u00e3o informado\"\n\n$d.Save()\n"}
